$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "31.064.11"
$ws.Range("E2").Value = "  +1.19%  "

# Row 3
$ws.Range("D3").Value = "1.954.50"
$ws.Range("E3").Value = "  +0.26%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.70%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.03%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4906"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.45%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2973"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.08%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06847"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.26%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.14"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.48%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "107.79"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.17%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.947.00"
$ws.Range("E12").Value = "  -0.10%  "

# Row 13
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07754"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.13%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.449"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.36%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7070"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.23%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "283.00"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.52%  "

# Row 17
$ws.Range("D17").Value = "31.087.59"
$ws.Range("E17").Value = "  +1.19%  "

# Row 18
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.33%  "

# Row 19
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007760"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.54%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9998"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.04%  "

# Row 21
$ws.Range("D21").Value = "2.191.24"
$ws.Range("E21").Value = "  -0.45%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.503"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.15%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.20%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.506"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.41%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.816"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.48%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.57%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.40%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.213"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.93%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1056"
$ws.Range("D29").Style = "Normal"

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.423"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.11%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.581"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.83%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.573"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.90%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.448"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.35%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04964"
$ws.Range("D34").Style = "Normal"

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7572"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.61%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.180"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.64%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.728"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.14%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02035"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.34%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.703"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.18%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.171"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.61%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.458"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +9.06%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "73.98"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.30%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4497"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.43%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "109.29"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.68%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8818"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.82%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.187"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +10.90%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9999"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.29%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.444"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.73%  "

# Row 49
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "961.04"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.83%  "

# Row 50
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1267"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.89%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.2584"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.70%  "
